$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 7.005616666666666
$ws.Range("H2").Value2 = 21.01685
$ws.Range("I2").Value2 = 0.3190354092819215
$ws.Range("J2").Value2 = 0.3190354092819216
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 58.88086666666666
$ws.Range("N2").Value2 = 176.6426
$ws.Range("O2").Value2 = 0.2818386429293818
$ws.Range("P2").Value2 = 0.2818386429293819
$ws.Range("Q2").Value2 = 412.4967808677777
$ws.Range("R2").Value2 = 3712.471027809999
$ws.Range("S2").Value2 = 0.08991650679843667
$ws.Range("T2").Value2 = 0.0899165067984367

$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 7.005616666666666
$ws.Range("H3").Value2 = 21.01685
$ws.Range("I3").Value2 = 0.3190354092819215
$ws.Range("J3").Value2 = 0.3190354092819216
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 110.7417803333333
$ws.Range("N3").Value2 = 332.225341
$ws.Range("O3").Value2 = 0.5300756400448766
$ws.Range("P3").Value2 = 0.5300756400448767
$ws.Range("Q3").Value2 = 775.8144619995388
$ws.Range("R3").Value2 = 6982.33015799585
$ws.Range("S3").Value2 = 0.1691128987720937
$ws.Range("T3").Value2 = 0.1691128987720938

$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 7.005616666666666
$ws.Range("H4").Value2 = 21.01685
$ws.Range("I4").Value2 = 0.3190354092819215
$ws.Range("J4").Value2 = 0.3190354092819216
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 39.29429233333334
$ws.Range("N4").Value2 = 117.882877
$ws.Range("O4").Value2 = 0.1880857170257415
$ws.Range("P4").Value2 = 0.1880857170257415
$ws.Range("Q4").Value2 = 275.2807492752722
$ws.Range("R4").Value2 = 2477.52674347745
$ws.Range("S4").Value2 = 0.06000600371139112
$ws.Range("T4").Value2 = 0.06000600371139113

$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 14.03763466666667
$ws.Range("H5").Value2 = 42.112904
$ws.Range("I5").Value2 = 0.6392731338754509
$ws.Range("J5").Value2 = 0.639273133875451
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 58.88086666666666
$ws.Range("N5").Value2 = 176.6426
$ws.Range("O5").Value2 = 0.2818386429293818
$ws.Range("P5").Value2 = 0.2818386429293819
$ws.Range("Q5").Value2 = 826.5480951233777
$ws.Range("R5").Value2 = 7438.932856110399
$ws.Range("S5").Value2 = 0.1801718725126701
$ws.Range("T5").Value2 = 0.1801718725126702

$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 14.03763466666667
$ws.Range("H6").Value2 = 42.112904
$ws.Range("I6").Value2 = 0.6392731338754509
$ws.Range("J6").Value2 = 0.639273133875451
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 110.7417803333333
$ws.Range("N6").Value2 = 332.225341
$ws.Range("O6").Value2 = 0.5300756400448766
$ws.Range("P6").Value2 = 0.5300756400448767
$ws.Range("Q6").Value2 = 1554.552654655585
$ws.Range("R6").Value2 = 13990.97389190026
$ws.Range("S6").Value2 = 0.3388631156025237
$ws.Range("T6").Value2 = 0.3388631156025239

$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 14.03763466666667
$ws.Range("H7").Value2 = 42.112904
$ws.Range("I7").Value2 = 0.6392731338754509
$ws.Range("J7").Value2 = 0.639273133875451
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 39.29429233333334
$ws.Range("N7").Value2 = 117.882877
$ws.Range("O7").Value2 = 0.1880857170257415
$ws.Range("P7").Value2 = 0.1880857170257415
$ws.Range("Q7").Value2 = 551.5989202605343
$ws.Range("R7").Value2 = 4964.390282344808
$ws.Range("S7").Value2 = 0.120238145760257
$ws.Range("T7").Value2 = 0.1202381457602571

$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 0.915492
$ws.Range("H8").Value2 = 2.746476
$ws.Range("I8").Value2 = 0.04169145684262745
$ws.Range("J8").Value2 = 0.04169145684262746
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 58.88086666666666
$ws.Range("N8").Value2 = 176.6426
$ws.Range("O8").Value2 = 0.2818386429293818
$ws.Range("P8").Value2 = 0.2818386429293819
$ws.Range("Q8").Value2 = 53.90496238639999
$ws.Range("R8").Value2 = 485.1446614775999
$ws.Range("S8").Value2 = 0.01175026361827501
$ws.Range("T8").Value2 = 0.01175026361827501

$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 0.915492
$ws.Range("H9").Value2 = 2.746476
$ws.Range("I9").Value2 = 0.04169145684262745
$ws.Range("J9").Value2 = 0.04169145684262746
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 110.7417803333333
$ws.Range("N9").Value2 = 332.225341
$ws.Range("O9").Value2 = 0.5300756400448766
$ws.Range("P9").Value2 = 0.5300756400448767
$ws.Range("Q9").Value2 = 101.383213960924
$ws.Range("R9").Value2 = 912.448925648316
$ws.Range("S9").Value2 = 0.0220996256702591
$ws.Range("T9").Value2 = 0.0220996256702591

$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 0.915492
$ws.Range("H10").Value2 = 2.746476
$ws.Range("I10").Value2 = 0.04169145684262745
$ws.Range("J10").Value2 = 0.04169145684262746
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 39.29429233333334
$ws.Range("N10").Value2 = 117.882877
$ws.Range("O10").Value2 = 0.1880857170257415
$ws.Range("P10").Value2 = 0.1880857170257415
$ws.Range("Q10").Value2 = 35.973610276828
$ws.Range("R10").Value2 = 323.762492491452
$ws.Range("S10").Value2 = 0.007841567554093341
$ws.Range("T10").Value2 = 0.007841567554093341

